$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking correct value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total marks (B12): 78 -> 130
$ws.Range("B12").Value = 130

# Update correct/total marks text (E12): "76/84" -> "130/140"
$ws.Range("E12").Value = "130/140"
